# Commit: "testing: xlsx report testing"
# - Rename the worksheet from "Sheet1" to "Estimates" (also updates the
#   _xlnm._FilterDatabase defined names that point at the sheet).
# - Add a new shared string "Max (P=95%)" and use it for cell A65, which
#   previously (incorrectly) duplicated the "Min (P=95%)" label from A64.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Estimates"

$ws.Range("A65").Value = "Max (P=95%)"
